$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refine the xClutch map (column D, rows 7-16) to use the full pedal-travel
# range instead of only half of it. These were previously partly formula
# driven (D7+75, D11+16, ...); replace with plain literal values so the
# series covers 1100-1900 evenly instead of 1500-1900.
$ws.Range("D7").Value = 1100
$ws.Range("D8").Value = 1180
$ws.Range("D9").Value = 1260
$ws.Range("D10").Value = 1340
$ws.Range("D11").Value = 1420
$ws.Range("D12").Value = 1500
$ws.Range("D13").Value = 1580
$ws.Range("D14").Value = 1660
$ws.Range("D15").Value = 1720
$ws.Range("D16").Value = 1800
# D17 stays 1900 (unchanged)

# Widen the chart's value-axis minimum so the fuller data range is visible.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.Axes(2).MinimumScale = 1000

# Move the active selection (as last left by the author) to F18.
$ws.Range("F18").Select()
